$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 15:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1036417
$ws.Range("C4").Value = 652
$ws.Range("D4").Value = 143098
$ws.Range("E4").Value = 834035
$ws.Range("F4").Value = 19098
$ws.Range("G4").Value = 18
$ws.Range("H4").Value = 59284

# Row 9 - Alemania
$ws.Range("B9").Value = 160479
$ws.Range("C9").Value = 567
$ws.Range("E9").Value = 33765

# Row 18 - India
$ws.Range("B18").Value = 31787
$ws.Range("C18").Value = 463
$ws.Range("D18").Value = 7796
$ws.Range("E18").Value = 22983

# Row 23 - Arabia Saudita
$ws.Range("F23").Value = 125

# Row 42 - Serbia
$ws.Range("F42").Value = 78

# Row 56 - Argentina
$ws.Range("D56").Value = 1192
$ws.Range("E56").Value = 2728

# Rows 69-71: reorder Armenia/Irak/Afganistan -> Afganistan/Armenia/Irak
$ws.Range("A69").Value = "Afganistan"
$ws.Range("B69").Value = 1939
$ws.Range("C69").Value = 111
$ws.Range("D69").Value = 252
$ws.Range("E69").Value = 1627
$ws.Range("F69").Value = 7
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 60

$ws.Range("A70").Value = "Armenia"
$ws.Range("B70").Value = 1932
$ws.Range("C70").Value = 65
$ws.Range("D70").Value = 900
$ws.Range("E70").Value = 1002
$ws.Range("F70").Value = 10
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 30

$ws.Range("A71").Value = "Irak"
$ws.Range("B71").Value = 1928
$ws.Range("C71").Value = 0
$ws.Range("D71").Value = 1319
$ws.Range("E71").Value = 519
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = 90

# Rows 80-82: reorder Bulgaria/Cuba/Republica de Macedonia -> Republica de Macedonia/Bulgaria/Cuba
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 1442
$ws.Range("C80").Value = 21
$ws.Range("D80").Value = 627
$ws.Range("E80").Value = 742
$ws.Range("F80").Value = 13
$ws.Range("G80").Value = 2
$ws.Range("H80").Value = 73

$ws.Range("A81").Value = "Bulgaria"
$ws.Range("B81").Value = 1437
$ws.Range("C81").Value = 38
$ws.Range("D81").Value = 243
$ws.Range("E81").Value = 1133
$ws.Range("F81").Value = 39
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 61

$ws.Range("A82").Value = "Cuba"
$ws.Range("B82").Value = 1437
$ws.Range("C82").Value = 0
$ws.Range("D82").Value = 575
$ws.Range("E82").Value = 804
$ws.Range("F82").Value = 12
$ws.Range("G82").Value = 0
$ws.Range("H82").Value = 58
